$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.060.81'
$ws.Range("E2").Value = '  -3.74%  '
$ws.Range("D3").Value = '1.640.90'
$ws.Range("E3").Value = '  -3.68%  '
$ws.Range("E4").Value = '  +0.35%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '307.48'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.76%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3911'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.15%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3853'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -4.60%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.32%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.352'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -8.19%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '49.04'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -8.28%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.08466'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -4.02%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '23.95'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -7.90%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.121'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.82%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.00001280'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -5.67%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '7.476'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -6.43%  '
$ws.Range("D17").Value = '1.651.24'
$ws.Range("E17").Value = '  -5.37%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '94.25'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.79%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06947'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.74%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '20.77'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.06%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.923'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -5.67%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.24%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '13.67'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -4.83%  '
$ws.Range("D24").Value = '24.070.22'
$ws.Range("E24").Value = '  -3.69%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.347'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.76%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.682'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -9.82%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '22.44'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -4.96%  '
$ws.Range("E28").Value = '  +5.46%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '158.33'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.85%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '141.84'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -5.99%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '5.320'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -13.17%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.468'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -6.09%  '
$ws.Range("D33").Value = '1.776.81'
$ws.Range("E33").Value = '  -10.84%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '7.136'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.68%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.08056'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -5.84%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.02927'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -7.79%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.9726'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -7.43%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.2702'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -6.40%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.09242'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.39%  '
$ws.Range("E40").Value = '  -1.17%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '10.00'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -10.09%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.7615'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -8.70%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '13.09'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -7.05%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '15.95'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -7.36%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.6880'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -7.13%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.480'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -8.09%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.083'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -4.15%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.08380'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -4.82%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '134.05'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.32%  '
$ws.Range("E51").Value = '  -10.58%  '
